# Updates the StateHospitalSpclty sheet data so that the State/Hospital
# pairing lines up with the correct Speciality/Organ rows, and refreshes
# the saved view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-by-row cell content corrections (State, Hospital, Speciality, Organ) ---

$ws.Range("A5").Value  = "Tamilnadu"
$ws.Range("B5").Value  = "RELA"
$ws.Range("C5").Value  = "Heart"
$ws.Range("D5").Value  = "Cardiology"

$ws.Range("A6").Value  = "Odisha"
$ws.Range("C6").Value  = "Lungs"
$ws.Range("D6").Value  = "Pulomonology"

$ws.Range("A7").Value  = "Odisha"
$ws.Range("C7").Value  = "Liver"
$ws.Range("D7").Value  = "Hepatology"

$ws.Range("A8").Value  = "Odisha"
$ws.Range("C8").Value  = "Kidney"
$ws.Range("D8").Value  = "Nephrology"

$ws.Range("A9").Value  = "Odisha"
$ws.Range("C9").Value  = "Bones"
$ws.Range("D9").Value  = "Orthopaedic"

$ws.Range("A10").Value = "Odisha"
$ws.Range("C10").Value = "Children"
$ws.Range("D10").Value = "Pediatrics"

$ws.Range("A11").Value = "Tamilnadu"
$ws.Range("B11").Value = "RELA"
$ws.Range("C11").Value = "Lungs"
$ws.Range("D11").Value = "Pulomonology"

$ws.Range("A12").Value = "Tamilnadu"
$ws.Range("B12").Value = "RELA"
$ws.Range("C12").Value = "Liver"
$ws.Range("D12").Value = "Hepatology"

$ws.Range("B13").Value = "AIMS"
$ws.Range("C13").Value = "Kidney"
$ws.Range("D13").Value = "Nephrology"

$ws.Range("B14").Value = "AIMS"
$ws.Range("C14").Value = "Bones"
$ws.Range("D14").Value = "Orthopaedic"

$ws.Range("A15").Value = "Delhi"
$ws.Range("B15").Value = "AIMS"
$ws.Range("C15").Value = "Children"
$ws.Range("D15").Value = "Pediatrics"

$ws.Range("A16").Value = "Delhi"
$ws.Range("B16").Value = "AIMS"
$ws.Range("C16").Value = "Eyes"
$ws.Range("D16").Value = "Opthalmology"

$ws.Range("A17").Value = "Kolkata"
$ws.Range("B17").Value = "Apollo"

$ws.Range("A18").Value = "Kolkata"
$ws.Range("B18").Value = "Apollo"

$ws.Range("A19").Value = "Kolkata"
$ws.Range("B19").Value = "Apollo"

$ws.Range("A20").Value = "Kolkata"
$ws.Range("B20").Value = "Apollo"

$ws.Range("A21").Value = "Kolkata"
$ws.Range("B21").Value = "Apollo"
$ws.Range("C21").Value = "Bone Health"
$ws.Range("D21").Value = "Rheumatology"

$ws.Range("A22").Value = "Delhi"
$ws.Range("B22").Value = "AIMS"

$ws.Range("A23").Value = "Delhi"
$ws.Range("B23").Value = "AIMS"

$ws.Range("A24").Value = "Delhi"
$ws.Range("B24").Value = "Medanta"

$ws.Range("A25").Value = "Delhi"
$ws.Range("B25").Value = "Medanta"

$ws.Range("A26").Value = "Delhi"
$ws.Range("B26").Value = "Medanta"

$ws.Range("A27").Value = "Delhi"
$ws.Range("B27").Value = "Medanta"

$ws.Range("A28").Value = "Tamilnadu"
$ws.Range("B28").Value = "RELA"
$ws.Range("C28").Value = "Kidney"
$ws.Range("D28").Value = "Nephrology"

$ws.Range("A29").Value = "Kolkata"
$ws.Range("B29").Value = "Woodlands"
$ws.Range("C29").Value = "Heart"
$ws.Range("D29").Value = "Cardiology"

$ws.Range("B30").Value = "Woodlands"
$ws.Range("C30").Value = "Lungs"
$ws.Range("D30").Value = "Pulomonology"

$ws.Range("B31").Value = "Woodlands"
$ws.Range("C31").Value = "Liver"
$ws.Range("D31").Value = "Hepatology"

$ws.Range("B32").Value = "Woodlands"
$ws.Range("C32").Value = "Kidney"
$ws.Range("D32").Value = "Nephrology"

$ws.Range("A33").Value = "Tamilnadu"
$ws.Range("B33").Value = "RELA"
$ws.Range("C33").Value = "Bones"
$ws.Range("D33").Value = "Orthopaedic"

$ws.Range("A34").Value = "Tamilnadu"
$ws.Range("B34").Value = "RELA"
$ws.Range("C34").Value = "Children"
$ws.Range("D34").Value = "Pediatrics"

$ws.Range("A35").Value = "Tamilnadu"
$ws.Range("B35").Value = "RELA"
$ws.Range("C35").Value = "Eyes"
$ws.Range("D35").Value = "Opthalmology"

$ws.Range("A36").Value = "Tamilnadu"
$ws.Range("B36").Value = "RELA"
$ws.Range("C36").Value = "Nerve"
$ws.Range("D36").Value = "Neurology"

$ws.Range("A37").Value = "Tamilnadu"
$ws.Range("B37").Value = "RELA"
$ws.Range("C37").Value = "Mental"
$ws.Range("D37").Value = "Psychology"

$ws.Range("A38").Value = "Odisha"
$ws.Range("B38").Value = "AIMS"
$ws.Range("C38").Value = "Heart"
$ws.Range("D38").Value = "Cardiology"

$ws.Range("B39").Value = "Apollo"

$ws.Range("B40").Value = "Apollo"

$ws.Range("B41").Value = "Disha"
$ws.Range("C41").Value = "Eyes"
$ws.Range("D41").Value = "Opthalmology"

$ws.Range("B42").Value = "Anandalok"
$ws.Range("C42").Value = "Eyes"
$ws.Range("D42").Value = "Opthalmology"

# --- Style updates ---
# Rows 29-32 pick up the explicit black-font style (matches cellXfs index 2)
$ws.Range("C29:D32").Font.Color = 0

# Rows 39-42 lose that explicit style, reverting to the sheet's default style
$ws.Range("C39:D42").ClearFormats()

# --- View / selection state ---
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H35").Select()
